$p = $ppt.ActivePresentation

# Slide 4 (sldId 262) notes body placeholder ("Notizenplatzhalter 2"):
#  - add "Automation" to the list of engineering backgrounds
#  - start a new paragraph at "experiences" and mention "Scrum, Java" skills
#  - drop the "Everybody will state his/her opinion and" lead-in sentence
$slide = $p.Slides.Item(4)
$notesShape = $slide.NotesPage.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

$text = $tr.Text
# PowerPoint reports paragraph breaks as CR (`r) when reading .Text back, but a
# bare CR does not start a new paragraph when written back through this COM
# layer - only LF (`n) does. Normalize so the logic below is symmetric.
$text = $text.Replace([char]13, [char]10)

# 1) "Industrial Engineering), " -> "Industrial Engineering, Automation), " and
#    split the paragraph right before "experiences".
$text = $text.Replace(
    " Engineering, Industrial Engineering), experiences",
    " Engineering, Industrial Engineering, Automation), `nexperiences"
)

# 2) "experiences (" -> "experiences (Scrum, Java, "
$text = $text.Replace("experiences (", "experiences (Scrum, Java, ")

# 3) Drop the "Everybody will state his/ her opinion and " lead-in before "we will support"
$text = $text.Replace("Everybody will state his/ her opinion and we will support", "we will support")

$tr.Text = $text
